# Apply updated dSF (column F) values as part of a data repull / push / mean
# recalculation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    6  = 0
    9  = -1
    11 = -10
    12 = -1
    17 = -2
    22 = -2
    23 = -1
    24 = -3
    25 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
